# Auto-generated edit script to apply the numeric / URL updates described in the diff.
# Sheet mapping: Worksheets.Item(1) = 展览, Item(2) = 演出, Item(3) = 本地生活, Item(4) = 全部类型
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet1) ---
$ws1.Range("F2").Value = 10048
$ws1.Range("F4").Value = 2520
$ws1.Range("F6").Value = 282
$ws1.Range("F9").Value = 754
$ws1.Range("F12").Value = 1035
$ws1.Range("F13").Value = 3117
$ws1.Range("F14").Value = 2343
$ws1.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202403/JmFXyFgc1710844373405.jpeg"
$ws1.Range("F16").Value = 2058
$ws1.Range("F20").Value = 1582
$ws1.Range("F21").Value = 546
$ws1.Range("F22").Value = 49
$ws1.Range("F23").Value = 230
$ws1.Range("F24").Value = 3
$ws1.Range("F25").Value = 8
$ws1.Range("F29").Value = 75
$ws1.Range("F30").Value = 355
$ws1.Range("F31").Value = 571
$ws1.Range("F33").Value = 215
$ws1.Range("F35").Value = 12
$ws1.Range("F36").Value = 303
$ws1.Range("F37").Value = 1639
$ws1.Range("F38").Value = 99
$ws1.Range("F39").Value = 407
$ws1.Range("F40").Value = 48
$ws1.Range("F41").Value = 432
$ws1.Range("F42").Value = 920
$ws1.Range("F44").Value = 344

# --- 演出 (sheet2) ---
$ws2.Range("F9").Value = 6

# --- 全部类型 (sheet4) ---
$ws4.Range("F2").Value = 10048
$ws4.Range("F4").Value = 2520
$ws4.Range("F8").Value = 282
$ws4.Range("F11").Value = 754
$ws4.Range("F13").Value = 1035
$ws4.Range("F14").Value = 3117
$ws4.Range("F15").Value = 2343
$ws4.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202403/JmFXyFgc1710844373405.jpeg"
$ws4.Range("F16").Value = 2058
$ws4.Range("F17").Value = 2058
$ws4.Range("F20").Value = 1582
$ws4.Range("F21").Value = 546
$ws4.Range("F22").Value = 49
$ws4.Range("F23").Value = 230
$ws4.Range("F24").Value = 3
$ws4.Range("F25").Value = 8
$ws4.Range("F29").Value = 75
$ws4.Range("F30").Value = 355
$ws4.Range("F31").Value = 571
$ws4.Range("F36").Value = 215
$ws4.Range("F38").Value = 12
$ws4.Range("F40").Value = 303
$ws4.Range("F41").Value = 1639
$ws4.Range("F42").Value = 99
$ws4.Range("F44").Value = 407
$ws4.Range("F45").Value = 48
$ws4.Range("F46").Value = 432
$ws4.Range("F47").Value = 921
$ws4.Range("F49").Value = 344
$ws4.Range("F50").Value = 6

Write-Host "Applied all updates."
